$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 195; this shifts existing rows 195-210 down to 196-211
$ws.Rows.Item(195).Insert()

# Populate the newly inserted row 195 with the new weekly price record
$ws.Range("A195").Value = 3
$ws.Range("B195").Value = "Femacal de La Calera"
$ws.Range("C195").Value = "Coquimbo"
$ws.Range("D195").Value = 44461
$ws.Range("D195").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E195").Value = 5
$ws.Range("F195").Value = 100112031
$ws.Range("G195").Value = "Poroto verde"
$ws.Range("H195").Value = "Magnum"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 35
$ws.Range("K195").Value = 33000
$ws.Range("L195").Value = 33000
$ws.Range("M195").Value = 33000
$ws.Range("N195").Value = "`$/malla 25 kilos"
$ws.Range("O195").Value = "Provincia de Quillota"
$ws.Range("P195").Value = 1320
$ws.Range("Q195").Value = 25
$ws.Range("R195").Value = "Hortaliza"
